$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 167 (pushes existing rows 167:223 down to 171:227),
# inheriting formatting (e.g. date style on column D) from the row above, as Excel does.
$ws.Rows("167:170").Insert()

# Populate the 4 newly inserted rows with the new weekly price entries.
# Columns A,B,C,E,F,G,H,I,J,R are constant for this product/market block,
# so copy them from the row directly below (the shifted former row167).
# Note: use Value2 when chaining a getter straight into a setter.
$constCols = @("A","B","C","E","F","G","H","I","J","R")
foreach ($col in $constCols) {
    $src = $ws.Range("$col" + "171").Value2
    for ($r = 167; $r -le 170; $r++) {
        $ws.Range("$col$r").Value2 = $src
    }
}

# Row 167: Artic Pride / Especial
$ws.Range("D167").Value = 44924
$ws.Range("K167").Value = "Artic Pride"
$ws.Range("L167").Value = "Especial"
$ws.Range("M167").Value = 16
$ws.Range("N167").Value = 450000
$ws.Range("O167").Value = 460000
$ws.Range("P167").Value = 455000
$ws.Range("Q167").Value = '$/bins (420 kilos)'
$ws.Range("S167").Value = 1083
$ws.Range("T167").Value = 420

# Row 168: Artic Pride / Primera
$ws.Range("D168").Value = 44924
$ws.Range("K168").Value = "Artic Pride"
$ws.Range("L168").Value = "Primera"
$ws.Range("M168").Value = 20
$ws.Range("N168").Value = 400000
$ws.Range("O168").Value = 410000
$ws.Range("P168").Value = 405000
$ws.Range("Q168").Value = '$/bins (420 kilos)'
$ws.Range("S168").Value = 964
$ws.Range("T168").Value = 420

# Row 169: Artic Pride / Segunda
$ws.Range("D169").Value = 44924
$ws.Range("K169").Value = "Artic Pride"
$ws.Range("L169").Value = "Segunda"
$ws.Range("M169").Value = 10
$ws.Range("N169").Value = 350000
$ws.Range("O169").Value = 360000
$ws.Range("P169").Value = 355000
$ws.Range("Q169").Value = '$/bins (420 kilos)'
$ws.Range("S169").Value = 845
$ws.Range("T169").Value = 420

# Row 170: Red Diamond / Primera
$ws.Range("D170").Value = 44924
$ws.Range("K170").Value = "Red Diamond"
$ws.Range("L170").Value = "Primera"
$ws.Range("M170").Value = 20
$ws.Range("N170").Value = 420000
$ws.Range("O170").Value = 430000
$ws.Range("P170").Value = 425000
$ws.Range("Q170").Value = '$/bins (420 kilos)'
$ws.Range("S170").Value = 1012
$ws.Range("T170").Value = 420
